$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "name" column (C) values: replace Harry Potter names with
# Lord of the Rings names (deletePet test data update)
$ws.Range("C2").Value = "Merry"
$ws.Range("C3").Value = "Pippin"
$ws.Range("C4").Value = "Aragorn"
$ws.Range("C5").Value = "Gimli"
$ws.Range("C6").Value = "Legolas"

# Move the active selection to C7
$ws.Range("C7").Select()
